# Removing discrimination and candidate generation phases from intro
#
# Slide 4: the "(5) Candidate discrimination" textbox (TextBox 12 / id=13)
# is renamed to "classification" and narrowed/shifted to the right, and
# the elbow connector glued to it (Elbow Connector 55 / id=56) is
# re-routed to match the textbox's new position.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

# --- TextBox 12 (id=13): "(5) Candidate discrimination" -> "classification"
$textBox = $s.Shapes.Item(5)

# Reposition/resize (y/height are unchanged by the edit, so only Left/Width
# are touched).
$textBox.Left = 688.8491821289062
$textBox.Width = 100.5020523071289

# Change just the second paragraph's run text ("discrimination" ->
# "classification"), leaving its run-level formatting (bold, size, etc.)
# untouched.
$textBox.TextFrame.TextRange.Paragraphs(2).Runs(1).Text = "classification"

# --- Elbow Connector 55 (id=56): re-route to match the textbox's new spot
$connector = $s.Shapes.Item(26)

$connector.Left = 606.5913696289062
$connector.Top = -14.005826950073242
$connector.Height = 265.01763916015625
